# Mise à jour du classement - 01.04.2025 à 22:00
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("leaderboard2")
$ws1.Range("D3").Value = 584
$ws1.Range("D4").Value = 523
$ws1.Range("D5").Value = 40
$ws1.Range("B13").Value = "Dernière update le 01.04.25 à 22:00"

$ws2 = $wb.Worksheets.Item("leaderboard3")
$ws2.Range("B13").Value = "Dernière update le 01.04.25 à 22:00"

$ws3 = $wb.Worksheets.Item("leaderboard4")
$ws3.Range("D4").Value = 9
$ws3.Range("D5").Value = 3
$ws3.Range("B13").Value = "Dernière update le 01.04.25 à 22:00"
